# NEBRASKA_2020.xlsx cleanup:
#  - rename header row to short machine-friendly column names
#  - title-case the Spanish connector words (de/del/el/la/los/las/y) inside
#    state/municipality names (but leave connectors that are already
#    capitalized, e.g. "El Porvenir", "La Paz", untouched)
#  - fix the one all-caps state name ("GUANAJUATO" -> "Guanajuato")
#  - drop the trailing footnote/source rows after the data table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function TitleCaseConnectors([string]$s) {
    # .NET regex is case-sensitive by default, so this only touches the
    # lower-case forms ("de", "del", "el", "la", "los", "las", "y") and
    # leaves already-capitalized words (e.g. "El", "La") untouched.
    $s = [regex]::Replace($s, '\bdel\b', 'Del')
    $s = [regex]::Replace($s, '\bde\b', 'De')
    $s = [regex]::Replace($s, '\blos\b', 'Los')
    $s = [regex]::Replace($s, '\blas\b', 'Las')
    $s = [regex]::Replace($s, '\bel\b', 'El')
    $s = [regex]::Replace($s, '\bla\b', 'La')
    $s = [regex]::Replace($s, '\by\b', 'Y')
    return $s
}

# 1) Rename header columns to machine-friendly names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Walk the data rows (2..631) and fix the state (A) / municipality (B)
#    text cells: title-case the lower-case connector words.
#    NOTE: this engine's `-eq`/`-ne` string comparisons are case-INsensitive,
#    so case-sensitive equality checks below use the .NET `.Equals()`
#    instance method (ordinal, case-sensitive) instead.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($colLetter in @("A", "B")) {
        $cell = $ws.Range($colLetter + $r)
        $v = $cell.Value2
        if ($v -ne $null -and $v -is [string] -and $v -ne "") {
            $nv = TitleCaseConnectors $v
            if (-not $nv.Equals($v)) {
                $cell.Value = $nv
            }
        }
    }
}

# 3) One-off fix: "GUANAJUATO" (all caps) -> "Guanajuato".
$a163 = $ws.Range("A163").Value2
if ($a163 -ne $null -and $a163.Equals("GUANAJUATO")) {
    $ws.Range("A163").Value = "Guanajuato"
}

# 4) Drop the trailing footnote/source rows (633-637) that sit below the
#    data table (row 632 is already a blank spacer row).
$ws.Rows("633:637").Delete()
